$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.169.28'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '3.276.48'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '583.07'
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("D6").Value = '184.91'
$ws.Range("E6").Value = '  +1.58%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  +1.16%  '

$ws.Range("E9").Value = '  -1.87%  '

$ws.Range("D10").Value = '6.59'
$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("E11").Value = '  -2.74%  '

$ws.Range("D12").Value = '3.847.23'
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("E14").Value = '  -3.24%  '

$ws.Range("D15").Value = '68.174.84'
$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("E16").Value = '  -1.34%  '

$ws.Range("D17").Value = '3.278.34'
$ws.Range("E17").Value = '  +3.53%  '

$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").Value = '416.35'
$ws.Range("E20").Value = '  +5.88%  '

$ws.Range("D21").Value = '7.57'
$ws.Range("E21").Value = '  -1.32%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = '71.50'
$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("E25").Value = '  -1.24%  '

$ws.Range("E26").Value = '  -2.06%  '

$ws.Range("D27").Value = '9.50'
$ws.Range("E27").Value = '  -0.89%  '

$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("E29").Value = '  -1.75%  '

$ws.Range("D30").Value = '22.78'
$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("E31").Value = '  -3.62%  '

$ws.Range("D32").Value = '6.92'
$ws.Range("E32").Value = '  -2.80%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.25'
$ws.Range("E34").Value = '  -2.69%  '

$ws.Range("D35").Value = '164.35'
$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("D38").Value = '27.39'
$ws.Range("E38").Value = '  +3.97%  '

$ws.Range("E39").Value = '  -3.23%  '

$ws.Range("D40").Value = '4.49'
$ws.Range("E40").Value = '  -1.98%  '

$ws.Range("D41").Value = '6.37'
$ws.Range("E41").Value = '  -3.72%  '

$ws.Range("D42").Value = '2.671.16'
$ws.Range("E42").Value = '  +2.62%  '

$ws.Range("D43").Value = '40.85'
$ws.Range("E43").Value = '  -1.23%  '

$ws.Range("E44").Value = '  -1.25%  '

$ws.Range("D45").Value = '0.0681'
$ws.Range("E45").Value = '  -1.20%  '

$ws.Range("D46").Value = '337.56'
$ws.Range("E46").Value = '  -2.37%  '

$ws.Range("D47").Value = '24.67'
$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").Value = '0.0275'
$ws.Range("E48").Value = '  -2.63%  '

$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("E50").Value = '  -1.03%  '

$ws.Range("D51").Value = '0.975'
$ws.Range("E51").Value = '  -0.73%  '
